# Update Department_Category helper sheet strings so that reference links
# (MC2 Link / Has Departments / Department_Category values) are refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B12: Backend_Support_Unit -> COE
$ws.Range("B12").Value = "COE"

# C2:C12 all hold the "LINK" label -> "Click Me"
$ws.Range("C2").Value = "Click Me"
$ws.Range("C3").Value = "Click Me"
$ws.Range("C4").Value = "Click Me"
$ws.Range("C5").Value = "Click Me"
$ws.Range("C6").Value = "Click Me"
$ws.Range("C7").Value = "Click Me"
$ws.Range("C8").Value = "Click Me"
$ws.Range("C9").Value = "Click Me"
$ws.Range("C10").Value = "Click Me"
$ws.Range("C11").Value = "Click Me"
$ws.Range("C12").Value = "Click Me"

# D2 (User_Development): Communication_&_PR -> Communication_and_PR
$ws.Range("D2").Value = "UserDevelopment_Team_Denmark_Company,UserDevelopment_Team_Netherlands_Company,UserDevelopment_Team_Sweden_Company,UserDevelopment_Team_Germany_Company,Public_Relationship_Team,MarComm_Team,Event,Content_Team,Europe_UserDevelopment,Communication_and_PR"

# D3 (Sales_Operation): Fleet_&_Business_Sales -> Fleet_and_Business_Sales
$ws.Range("D3").Value = "UserTeam_Norway_Company,UserOperations_Team_Germany_Company,UserOperations_Team_Netherlands_Company,UserOperations_Team_Denmark_Company,UserOperations_Team_Sweden_Company,Fleet_Planning_Team,Fleet_Operation_Team,Fleet_Management_Team_Netherlands_Company,Fleet_Management_Team_Sweden_Company,Fleet_Management_Team_Denmark_Company,Business_Development,Business_Intelligence,KA_Management_Team,Commercial_Product,Europe_Commercial_Operation,Europe_UserOperation_Department,Europe_Business_Operation,Fleet_and_Business_Sales,Partner_Strategy,Retail_Sales,Sales_Planning"

# D6 (Power_Operation): "&" -> "and" in several entries, last entry renamed/underscored
$ws.Range("D6").Value = "Power_Business_Operation_Team,Power_Management_Team_Germany_Company,Power_Management_Team_Netherlands_Company,Strategy and Business_Development Team,Market_Launch and Enabling Team,Power_Management_Team_Sweden_Company,Power_Management_Team_Denmark_Company,Europe_Power_Operation_Department,Power_Operation,Power_market_launch_and_enabling_team"

# D7 (Service_Operation): Parts_&_Logistics -> Parts_and_Logistics
$ws.Range("D7").Value = "Service_PMO,Service_Planning_Team,Service_Quality_Team,Spare_Parts_Team,Service_Operations_Team_Germany_Company,Service_Team_Norway_Company,Service_Operations_Team_Denmark_Company,Service_Operations_Team_Netherlands_Company,Service_Operations_Team_Sweden_Company,Europe_Service_Operation_Department,Operation_Support,Parts_and_Logistics"

# D8 (Market_Entry_and_Operation): Europe_Project_Management_&_Enabling_Department -> ..._and_...
$ws.Range("D8").Value = "Market_Entry_PMO_Team,Europe_Project_Management_and_Enabling_Department,Regional_Operations_Support_Team ,Market_Expansion_Planning"

# D10 (System_Development): trailing space removed
$ws.Range("D10").Value = "Digital_PMO_Team"

# D12 (Backend_Support_Unit / COE row): Has Departments list fully replaced
$ws.Range("D12").Value = "Europe_Business_HRBP_Department,Controlling_and_Planning_Department,Legal_EU_Department,NIO_Life_Supply_Chain_Department,EHS,Digital_Development_PMO_Team,Product_Marketing_Department,Europe_Product_Experience_Department,Purchasing_Governance_and_BP_Team"
